# Apply updated confidence/frequency values and reordered category
# combinations as described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 5198

# Row 3
$ws.Range("C3").Value = 991

# Row 4
$ws.Range("C4").Value = 552

# Row 5
$ws.Range("A5").Value = "Fashion & Accessories"
$ws.Range("B5").Value = "Kitchen & Dining"
$ws.Range("C5").Value = 236

# Row 6
$ws.Range("A6").Value = "Textiles & Cozy Items"
$ws.Range("B6").Value = "Fashion & Accessories"
$ws.Range("C6").Value = 223

# Row 7
$ws.Range("B7").Value = "Storage & Organization"
$ws.Range("C7").Value = 211

# Row 8
$ws.Range("B8").Value = "Textiles & Cozy Items"
$ws.Range("C8").Value = 206

# Row 9
$ws.Range("B9").Value = "Kitchen & Dining"
$ws.Range("C9").Value = 200

# Row 10
$ws.Range("C10").Value = 197

# Row 11
$ws.Range("B11").Value = "Storage & Organization"
$ws.Range("C11").Value = 175
